$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number + new text values for the changed columns.
# Columns D (Price), E (Volume 1h), F (Data/date), G (Hora/hour) were
# refreshed by the GitHub Actions symbol-list updater on 2023-01-19.
$updates = @(
    @{Row=2; D="288.30"; E="-3.77%"; F="19-1-2023"; G="0"},
    @{Row=3; D="30.64"; E="-4.55%"; F="19-1-2023"; G="0"},
    @{Row=4; D="4.901"; E="-1.48%"; F="19-1-2023"; G="0"},
    @{Row=5; D="0.07094"; E="-9.84%"; F="19-1-2023"; G="0"},
    @{Row=6; D="1.782"; E="-19.73%"; F="19-1-2023"; G="0"},
    @{Row=7; D="7.679"; E="-1.67%"; F="19-1-2023"; G="0"},
    @{Row=8; D="3.776"; E="-1.79%"; F="19-1-2023"; G="0"},
    @{Row=9; D="0.8972"; E="-2.55%"; F="19-1-2023"; G="0"},
    @{Row=10; D="0.1634"; E="-5.99%"; F="19-1-2023"; G="0"},
    @{Row=11; D="0.07569"; E="-4.69%"; F="19-1-2023"; G="0"},
    @{Row=12; D="0.07943"; E="-7.82%"; F="19-1-2023"; G="0"},
    @{Row=13; D="0.03019"; E="-2.23%"; F="19-1-2023"; G="0"},
    @{Row=14; D="0.09972"; E="-0.31%"; F="19-1-2023"; G="0"},
    @{Row=15; D="0.001498"; E="-1.24%"; F="19-1-2023"; G="0"},
    @{Row=16; D="0.005712"; E="-2.62%"; F="19-1-2023"; G="0"},
    @{Row=17; F="19-1-2023"; G="0"},
    @{Row=18; D="3.466"; E="0.17%"; F="19-1-2023"; G="0"},
    @{Row=19; D="2.099"; E="-3.63%"; F="19-1-2023"; G="0"},
    @{Row=20; E="0.01%"; F="19-1-2023"; G="0"},
    @{Row=21; D="0.1298"; E="1.34%"; F="19-1-2023"; G="0"},
    @{Row=22; D="4.048"; E="-5.70%"; F="19-1-2023"; G="0"},
    @{Row=23; D="0.2002"; E="11.57%"; F="19-1-2023"; G="0"},
    @{Row=24; D="0.04503"; E="-1.99%"; F="19-1-2023"; G="0"},
    @{Row=25; E="-1.09%"; F="19-1-2023"; G="0"},
    @{Row=26; D="0.004633"; E="4.82%"; F="19-1-2023"; G="0"},
    @{Row=27; D="0.0001251"; E="0.15%"; F="19-1-2023"; G="0"},
    @{Row=28; F="19-1-2023"; G="0"},
    @{Row=29; F="19-1-2023"; G="0"},
    @{Row=30; F="19-1-2023"; G="0"},
    @{Row=31; F="19-1-2023"; G="0"},
    @{Row=32; F="19-1-2023"; G="0"},
    @{Row=33; F="19-1-2023"; G="0"},
    @{Row=34; F="19-1-2023"; G="0"},
    @{Row=35; F="19-1-2023"; G="0"},
    @{Row=36; F="19-1-2023"; G="0"},
    @{Row=37; F="19-1-2023"; G="0"},
    @{Row=38; F="19-1-2023"; G="0"},
    @{Row=39; E="-7.88%"; F="19-1-2023"; G="0"},
    @{Row=40; D="0.04341"; E="-8.71%"; F="19-1-2023"; G="0"},
    @{Row=41; D="0.007375"; E="-1.13%"; F="19-1-2023"; G="0"},
    @{Row=42; E="-3.73%"; F="19-1-2023"; G="0"},
    @{Row=43; D="0.002003"; E="-14.37%"; F="19-1-2023"; G="0"},
    @{Row=44; D="0.009304"; E="-20.65%"; F="19-1-2023"; G="0"},
    @{Row=45; D="0.00005872"; E="-1.61%"; F="19-1-2023"; G="0"},
    @{Row=46; D="0.00000000751"; E="0.15%"; F="19-1-2023"; G="0"},
    @{Row=47; D="2.227"; E="171.36%"; F="19-1-2023"; G="0"},
    @{Row=48; D="0.003004"; E="-11.30%"; F="19-1-2023"; G="0"},
    @{Row=49; D="0.00002102"; E="0.15%"; F="19-1-2023"; G="0"},
    @{Row=50; D="0.0002002"; E="0.15%"; F="19-1-2023"; G="0"},
    @{Row=51; F="19-1-2023"; G="0"}
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
    if ($u.ContainsKey("F")) {
        $cell = $ws.Cells.Item($u.Row, 6)
        $cell.NumberFormat = "@"
        $cell.Value = $u.F
    }
    if ($u.ContainsKey("G")) {
        $cell = $ws.Cells.Item($u.Row, 7)
        $cell.NumberFormat = "@"
        $cell.Value = $u.G
    }
}
